$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 5-10 (the extra backbone/source-well combinations) so only
# rows 2-4 (one row per reagent) remain.
$ws.Range("A5:I10").EntireRow.Delete()

# Update the remaining rows to reflect the new source wells / destination
# wells / volumes / reagents.
# Row 3 (UID 2): source well A2, destination well A1, volume 125, reagent "DNA ligase"
$ws.Range("D3").Value = "A2"
$ws.Range("G3").Value = "A1"
$ws.Range("H3").Value = 125
$ws.Range("I3").Value = "DNA ligase"

# Row 4 (UID 3): source well A3, destination well A1, volume 250, reagent "BsmBI (NEB)"
$ws.Range("D4").Value = "A3"
$ws.Range("G4").Value = "A1"
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = "BsmBI (NEB)"
